# Add the missing "9:00 - 10:00" interval row to the schedule.
# It belongs right after the "8:00 - 9:00" row (row 3), so insert a new
# row there and push the following rows down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("3:3").Insert()
$ws.Range("A3").Value = "9:00 - 10:00"

# Move the active selection to the newly inserted cell, matching the
# saved workbook state.
$ws.Range("A3").Select()
